$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 131 ----
$c = $ws.Range("A131")
$c.Formula = "testPrivacyPolicyAgreementsEmptyVersion"

$ws.Range("A1").Copy()
$ws.Range("B131").PasteSpecial(-4122)
$c = $ws.Range("B131")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "No"

$ws.Range("C2").Copy()
$ws.Range("C131").PasteSpecial(-4122)
$c = $ws.Range("C131")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "testdata-admin.xlsx,feeStructure"

$ws.Range("C2").Copy()
$ws.Range("D131").PasteSpecial(-4122)
$c = $ws.Range("D131")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "RunRangeOfIterations"

$ws.Range("E2").Copy()
$ws.Range("E131").PasteSpecial(-4122)
$c = $ws.Range("E131")
$c.Borders.Item(9).LineStyle = -4142
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F131").PasteSpecial(-4122)
$c = $ws.Range("F131")
$c.Borders.Item(9).LineStyle = -4142
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G131").PasteSpecial(-4122)
$c = $ws.Range("G131")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "System Settings - feeStructure"

$ws.Range("I34").Copy()
$ws.Range("H131").PasteSpecial(-4122)
$c = $ws.Range("H131")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-pmessage,`n-psecurityKey"

$ws.Range("I62").Copy()
$ws.Range("I131").PasteSpecial(-4122)
$c = $ws.Range("I131")
$c.Borders.Item(9).LineStyle = -4142
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestPrivacyPolicyAgreementsEmptyVersion,`n-perrMessage,`n-pcolor,`n-pelementName,`n-pIteration Num"

$ws.Rows.Item(131).RowHeight = 165

# ---- Row 132 ----
$ws.Range("J1").Copy()
$ws.Range("A132").PasteSpecial(-4122)
$c = $ws.Range("A132")
$c.Formula = "testPrivacyPolicyAgreementsReschedulingwithSameDate"

$ws.Range("A1").Copy()
$ws.Range("B132").PasteSpecial(-4122)
$c = $ws.Range("B132")
$c.Value = "No"

$ws.Range("C2").Copy()
$ws.Range("C132").PasteSpecial(-4122)
$c = $ws.Range("C132")
$c.Value = "testdata-admin.xlsx,feeStructure"

$ws.Range("C2").Copy()
$ws.Range("D132").PasteSpecial(-4122)
$c = $ws.Range("D132")
$c.Value = "RunRangeOfIterations"

$ws.Range("E2").Copy()
$ws.Range("E132").PasteSpecial(-4122)
$c = $ws.Range("E132")
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F132").PasteSpecial(-4122)
$c = $ws.Range("F132")
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G132").PasteSpecial(-4122)
$c = $ws.Range("G132")
$c.Value = "System Settings - feeStructure"

$ws.Range("H2").Copy()
$ws.Range("H132").PasteSpecial(-4122)
$c = $ws.Range("H132")
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-pmessage,`n-psecurityKey"

$ws.Range("I62").Copy()
$ws.Range("I132").PasteSpecial(-4122)
$c = $ws.Range("I132")
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestPrivacyPolicyAgreementsReschedulingwithSameDate,`n-pAgreementList,`n-pfolderName,`n-pfileName,`n-ptosHeading,`n-pTitle,`n-pMessage"

$ws.Rows.Item(132).RowHeight = 165

# ---- Row 133 ----
$ws.Range("A1").Copy()
$ws.Range("A133").PasteSpecial(-4122)
$c = $ws.Range("A133")
$c.Formula = "testTermOfServiceAgreementsEmptyVersion"

$ws.Range("A1").Copy()
$ws.Range("B133").PasteSpecial(-4122)
$c = $ws.Range("B133")
$c.Formula = "no"

$ws.Range("C2").Copy()
$ws.Range("C133").PasteSpecial(-4122)
$c = $ws.Range("C133")
$c.Value = "testdata-admin.xlsx,feeStructure"

$ws.Range("C2").Copy()
$ws.Range("D133").PasteSpecial(-4122)
$c = $ws.Range("D133")
$c.Value = "RunRangeOfIterations"

$ws.Range("E2").Copy()
$ws.Range("E133").PasteSpecial(-4122)
$c = $ws.Range("E133")
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F133").PasteSpecial(-4122)
$c = $ws.Range("F133")
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G133").PasteSpecial(-4122)
$c = $ws.Range("G133")
$c.Value = "System Settings - feeStructure"

$ws.Range("I34").Copy()
$ws.Range("H133").PasteSpecial(-4122)
$c = $ws.Range("H133")
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-pmessage,`n-psecurityKey"

$ws.Range("I62").Copy()
$ws.Range("I133").PasteSpecial(-4122)
$c = $ws.Range("I133")
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestTermOfServiceAgreementsEmptyVersion,`n-perrMessage,`n-pcolor,`n-pelementName,`n-pIteration Num"

$ws.Rows.Item(133).RowHeight = 165

# ---- Row 134 ----
$ws.Range("J1").Copy()
$ws.Range("A134").PasteSpecial(-4122)
$c = $ws.Range("A134")
$c.Formula = "testTermOfServiceAgreementsReschedulingwithSameDate"

$ws.Range("A1").Copy()
$ws.Range("B134").PasteSpecial(-4122)
$c = $ws.Range("B134")
$c.Value = "No"

$ws.Range("C2").Copy()
$ws.Range("C134").PasteSpecial(-4122)
$c = $ws.Range("C134")
$c.Value = "testdata-admin.xlsx,feeStructure"

$ws.Range("C2").Copy()
$ws.Range("D134").PasteSpecial(-4122)
$c = $ws.Range("D134")
$c.Value = "RunRangeOfIterations"

$ws.Range("E2").Copy()
$ws.Range("E134").PasteSpecial(-4122)
$c = $ws.Range("E134")
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F134").PasteSpecial(-4122)
$c = $ws.Range("F134")
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G134").PasteSpecial(-4122)
$c = $ws.Range("G134")
$c.Value = "System Settings - feeStructure"

$ws.Range("H2").Copy()
$ws.Range("H134").PasteSpecial(-4122)
$c = $ws.Range("H134")
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-pmessage,`n-psecurityKey"

$ws.Range("I62").Copy()
$ws.Range("I134").PasteSpecial(-4122)
$c = $ws.Range("I134")
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestTermOfServiceAgreementsReschedulingwithSameDate,`n-pAgreementList,`n-pfolderName,`n-pfileName,`n-ptosHeading,`n-pTitle,`n-pMessage"

$ws.Rows.Item(134).RowHeight = 165

# ---- Row 135 ----
$c = $ws.Range("A135")
$c.Formula = "testApplicationAcknowledgementAgreementsEmptyVersion"

$ws.Range("A1").Copy()
$ws.Range("B135").PasteSpecial(-4122)
$c = $ws.Range("B135")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "No"

$ws.Range("C2").Copy()
$ws.Range("C135").PasteSpecial(-4122)
$c = $ws.Range("C135")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "testdata-admin.xlsx,feeStructure"

$ws.Range("C2").Copy()
$ws.Range("D135").PasteSpecial(-4122)
$c = $ws.Range("D135")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "RunRangeOfIterations"

$ws.Range("E2").Copy()
$ws.Range("E135").PasteSpecial(-4122)
$c = $ws.Range("E135")
$c.Borders.Item(9).LineStyle = -4142
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F135").PasteSpecial(-4122)
$c = $ws.Range("F135")
$c.Borders.Item(9).LineStyle = -4142
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G135").PasteSpecial(-4122)
$c = $ws.Range("G135")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "System Settings - feeStructure"

$ws.Range("I34").Copy()
$ws.Range("H135").PasteSpecial(-4122)
$c = $ws.Range("H135")
$c.Borders.Item(9).LineStyle = -4142
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-pmessage,`n-psecurityKey"

$ws.Range("I62").Copy()
$ws.Range("I135").PasteSpecial(-4122)
$c = $ws.Range("I135")
$c.Borders.Item(9).LineStyle = -4142
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestApplicationAcknowledgementAgreementsEmptyVersion,`n-perrMessage,`n-pcolor,`n-pelementName,`n-pIteration Num"

$ws.Rows.Item(135).RowHeight = 165

# ---- Row 136 ----
$ws.Range("J1").Copy()
$ws.Range("A136").PasteSpecial(-4122)
$c = $ws.Range("A136")
$c.Formula = "testApplicationAcknowledgementAgreementsReschedulingwithSameDate"

$ws.Range("A1").Copy()
$ws.Range("B136").PasteSpecial(-4122)
$c = $ws.Range("B136")
$c.Value = "No"

$ws.Range("C2").Copy()
$ws.Range("C136").PasteSpecial(-4122)
$c = $ws.Range("C136")
$c.Value = "testdata-admin.xlsx,feeStructure"

$ws.Range("C2").Copy()
$ws.Range("D136").PasteSpecial(-4122)
$c = $ws.Range("D136")
$c.Value = "RunRangeOfIterations"

$ws.Range("E2").Copy()
$ws.Range("E136").PasteSpecial(-4122)
$c = $ws.Range("E136")
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F136").PasteSpecial(-4122)
$c = $ws.Range("F136")
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G136").PasteSpecial(-4122)
$c = $ws.Range("G136")
$c.Value = "System Settings - feeStructure"

$ws.Range("H2").Copy()
$ws.Range("H136").PasteSpecial(-4122)
$c = $ws.Range("H136")
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-pmessage,`n-psecurityKey"

$ws.Range("I62").Copy()
$ws.Range("I136").PasteSpecial(-4122)
$c = $ws.Range("I136")
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestApplicationAcknowledgementAgreementsReschedulingwithSameDate,`n-pAgreementList,`n-pfolderName,`n-pfileName,`n-ptosHeading,`n-pTitle,`n-pMessage"

$ws.Rows.Item(136).RowHeight = 165

# ---- Row 137 ----
$ws.Range("J1").Copy()
$ws.Range("A137").PasteSpecial(-4122)
$c = $ws.Range("A137")
$c.Formula = "testEditPersonalFeeStructureReschedulingwithSameDate"

$ws.Range("A1").Copy()
$ws.Range("B137").PasteSpecial(-4122)
$c = $ws.Range("B137")
$c.Value = "No"

$ws.Range("A33").Copy()
$ws.Range("C137").PasteSpecial(-4122)
$c = $ws.Range("C137")
$c.Value = "testdata-admin.xlsx,feeStructure"

$ws.Range("C2").Copy()
$ws.Range("D137").PasteSpecial(-4122)
$c = $ws.Range("D137")
$c.Value = "RunOneIteration"

$ws.Range("E2").Copy()
$ws.Range("E137").PasteSpecial(-4122)
$c = $ws.Range("E137")
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F137").PasteSpecial(-4122)
$c = $ws.Range("F137")
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G137").PasteSpecial(-4122)
$c = $ws.Range("G137")
$c.Value = "System Settings - feeStructure"

$ws.Range("H2").Copy()
$ws.Range("H137").PasteSpecial(-4122)
$c = $ws.Range("H137")
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-psecurityKey"

$ws.Range("H2").Copy()
$ws.Range("I137").PasteSpecial(-4122)
$c = $ws.Range("I137")
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestEditPersonalFeeStructureReschedulingwithSameDate,`n-pAgreementList,`n-pexpEditHeading,`n-pamount,`n-ppercentage,`n-ptosHeading"

$ws.Rows.Item(137).RowHeight = 150

# ---- Row 138 ----
$ws.Range("J1").Copy()
$ws.Range("A138").PasteSpecial(-4122)
$c = $ws.Range("A138")
$c.Formula = "testEditMerchantFeeStructureReschedulingwithSameDate"

$ws.Range("A1").Copy()
$ws.Range("B138").PasteSpecial(-4122)
$c = $ws.Range("B138")
$c.Value = "No"

$ws.Range("A33").Copy()
$ws.Range("C138").PasteSpecial(-4122)
$c = $ws.Range("C138")
$c.Value = "testdata-admin.xlsx,feeStructure"

$ws.Range("C2").Copy()
$ws.Range("D138").PasteSpecial(-4122)
$c = $ws.Range("D138")
$c.Value = "RunOneIteration"

$ws.Range("E2").Copy()
$ws.Range("E138").PasteSpecial(-4122)
$c = $ws.Range("E138")
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F138").PasteSpecial(-4122)
$c = $ws.Range("F138")
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G138").PasteSpecial(-4122)
$c = $ws.Range("G138")
$c.Value = "System Settings - feeStructure"

$ws.Range("H2").Copy()
$ws.Range("H138").PasteSpecial(-4122)
$c = $ws.Range("H138")
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-psecurityKey"

$ws.Range("H2").Copy()
$ws.Range("I138").PasteSpecial(-4122)
$c = $ws.Range("I138")
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestEditMerchantFeeStructureReschedulingwithSameDate,`n-pAgreementList,`n-pexpEditHeading,`n-pamount,`n-ppercentage,`n-ptosHeading"

$ws.Rows.Item(138).RowHeight = 150

# ---- Row 139 ----
$ws.Range("J1").Copy()
$ws.Range("A139").PasteSpecial(-4122)
$c = $ws.Range("A139")
$c.Formula = "testViewPersonalAccountLimitsReschedulingwithSameDate"

$ws.Range("A1").Copy()
$ws.Range("B139").PasteSpecial(-4122)
$c = $ws.Range("B139")
$c.Value = "No"

$ws.Range("C2").Copy()
$ws.Range("C139").PasteSpecial(-4122)
$c = $ws.Range("C139")
$c.Value = "testdata-admin.xlsx,AccountLimits"

$ws.Range("C2").Copy()
$ws.Range("D139").PasteSpecial(-4122)
$c = $ws.Range("D139")
$c.Value = "RunOneIteration"

$ws.Range("E2").Copy()
$ws.Range("E139").PasteSpecial(-4122)
$c = $ws.Range("E139")
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F139").PasteSpecial(-4122)
$c = $ws.Range("F139")
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G139").PasteSpecial(-4122)
$c = $ws.Range("G139")
$c.Value = "SystemSettings -AccountLimit"

$ws.Range("H2").Copy()
$ws.Range("H139").PasteSpecial(-4122)
$c = $ws.Range("H139")
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-psecurityKey"

$ws.Range("H2").Copy()
$ws.Range("I139").PasteSpecial(-4122)
$c = $ws.Range("I139")
$c.Borders.Item(9).LineStyle = -4142
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestViewPersonalAccountLimitsReschedulingwithSameDate,`n-pAgreementList,`n-pexpEditHeading,`n-pamount,`n-ppercentage,`n-ptosHeading"

$ws.Rows.Item(139).RowHeight = 150

# ---- Row 140 ----
$ws.Range("J1").Copy()
$ws.Range("A140").PasteSpecial(-4122)
$c = $ws.Range("A140")
$c.Formula = "testViewMerchantAccountLimitsReschedulingwithSameDate"

$ws.Range("A1").Copy()
$ws.Range("B140").PasteSpecial(-4122)
$c = $ws.Range("B140")
$c.Value = "No"

$ws.Range("C2").Copy()
$ws.Range("C140").PasteSpecial(-4122)
$c = $ws.Range("C140")
$c.Value = "testdata-admin.xlsx,AccountLimits"

$ws.Range("C2").Copy()
$ws.Range("D140").PasteSpecial(-4122)
$c = $ws.Range("D140")
$c.Value = "RunOneIteration"

$ws.Range("E2").Copy()
$ws.Range("E140").PasteSpecial(-4122)
$c = $ws.Range("E140")
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F140").PasteSpecial(-4122)
$c = $ws.Range("F140")
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G140").PasteSpecial(-4122)
$c = $ws.Range("G140")
$c.Value = "SystemSettings -AccountLimit"

$ws.Range("H2").Copy()
$ws.Range("H140").PasteSpecial(-4122)
$c = $ws.Range("H140")
$c.Borders.Item(10).LineStyle = -4142
$c.Value = "coyni.admin.tests.LoginTest,`ntestAdminLogin,`n-ploginHeading,`n-ploginDescription,`n-pemail,`n-ppassword,`n-pauthyHeading,`n-pauthyDescription,`n-pcode,`n-psecurityKey"

$ws.Range("H2").Copy()
$ws.Range("I140").PasteSpecial(-4122)
$c = $ws.Range("I140")
$c.Formula = "coyni.admin.tests.SystemSettingsTest,`ntestViewMerchantAccountLimitsReschedulingwithSameDate,`n-pAgreementList,`n-pexpEditHeading,`n-pamount,`n-ppercentage,`n-ptosHeading"

$ws.Rows.Item(140).RowHeight = 150

# ---- Row 141 ----
$ws.Range("A1").Copy()
$ws.Range("A141").PasteSpecial(-4122)
$c = $ws.Range("A141")
$c.Value = "verify New Verification CodeSend Upto 5Times"

$ws.Range("A1").Copy()
$ws.Range("B141").PasteSpecial(-4122)
$c = $ws.Range("B141")
$c.Formula = "Yes"

$ws.Range("C2").Copy()
$ws.Range("C141").PasteSpecial(-4122)
$c = $ws.Range("C141")
$c.Value = "testdata-admin.xlsx,Login"

$ws.Range("C2").Copy()
$ws.Range("D141").PasteSpecial(-4122)
$c = $ws.Range("D141")
$c.Value = "RunRangeOfIterations"

$ws.Range("E2").Copy()
$ws.Range("E141").PasteSpecial(-4122)
$c = $ws.Range("E141")
$c.Formula = "'1"

$ws.Range("E2").Copy()
$ws.Range("F141").PasteSpecial(-4122)
$c = $ws.Range("F141")
$c.Formula = "'1"

$ws.Range("C2").Copy()
$ws.Range("G141").PasteSpecial(-4122)
$c = $ws.Range("G141")
$c.Value = "Login"

$ws.Range("H2").Copy()
$ws.Range("H141").PasteSpecial(-4122)
$c = $ws.Range("H141")
$c.Borders.Item(10).LineStyle = -4142
$c.Value = "coyni.admin.tests.LoginTest,`ntestNewCodeUpto5Times,`n-pemail,`n-ppassword,`n-pfiveTimeContent"

$ws.Range("I2").Copy()
$ws.Range("I141").PasteSpecial(-4122)
$c = $ws.Range("I141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("J141").PasteSpecial(-4122)
$c = $ws.Range("J141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("K141").PasteSpecial(-4122)
$c = $ws.Range("K141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("L141").PasteSpecial(-4122)
$c = $ws.Range("L141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("M141").PasteSpecial(-4122)
$c = $ws.Range("M141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("N141").PasteSpecial(-4122)
$c = $ws.Range("N141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("O141").PasteSpecial(-4122)
$c = $ws.Range("O141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("P141").PasteSpecial(-4122)
$c = $ws.Range("P141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("Q141").PasteSpecial(-4122)
$c = $ws.Range("Q141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("R141").PasteSpecial(-4122)
$c = $ws.Range("R141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("S141").PasteSpecial(-4122)
$c = $ws.Range("S141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("T141").PasteSpecial(-4122)
$c = $ws.Range("T141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("U141").PasteSpecial(-4122)
$c = $ws.Range("U141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("V141").PasteSpecial(-4122)
$c = $ws.Range("V141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("W141").PasteSpecial(-4122)
$c = $ws.Range("W141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("X141").PasteSpecial(-4122)
$c = $ws.Range("X141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("Y141").PasteSpecial(-4122)
$c = $ws.Range("Y141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("Z141").PasteSpecial(-4122)
$c = $ws.Range("Z141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("AA141").PasteSpecial(-4122)
$c = $ws.Range("AA141")
# (leave blank)

$ws.Range("O2").Copy()
$ws.Range("AB141").PasteSpecial(-4122)
$c = $ws.Range("AB141")
# (leave blank)

$ws.Rows.Item(141).RowHeight = 182.25
